$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells we touch stay text-formatted so that
# values such as "54.10" or "0.999" keep their exact printed form
# instead of being coerced into floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.781.90"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.524.97"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.42"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.19"
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.514.13"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  -2.74%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.198"
$ws.Range("E10").Value = "  +6.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.644"
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.10"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("E13").Value = "  -2.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.47"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.085.41"
$ws.Range("E15").Value = "  -2.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.33"
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.765.74"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.507.00"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.37"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "543.53"
$ws.Range("E21").Value = "  +11.12%  "
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.85"
$ws.Range("E23").Value = "  -8.27%  "
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.84"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.58"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.16"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.08"
$ws.Range("E29").Value = "  -3.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.09"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("E31").Value = "  -4.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.42"
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.33"
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "545.25"
$ws.Range("E35").Value = "  -4.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.411"
$ws.Range("E36").Value = "  +3.89%  "
$ws.Range("E37").Value = "  +4.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.10"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0761"
$ws.Range("E40").Value = "  -5.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.369.67"
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("E43").Value = "  -3.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.09"
$ws.Range("E44").Value = "  -6.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.52"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0439"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("E48").Value = "  -2.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("E49").Value = "  -6.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "136.85"
$ws.Range("E51").Value = "  +1.71%  "
